# Daily attendance processing - 2025-10-27 10:25:19
# Reorders the "Recorded By" (column G) entries on the "Session Analysis
# Results" sheet so that the literal "System" token is no longer first in
# the comma-separated list (the other recorder(s) come first instead).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Rows where "System, system, backup@backdoor.com" becomes
# "System, backup@backdoor.com, system"
$rowsTripleSystem = @(2, 29, 56)
foreach ($r in $rowsTripleSystem) {
    $ws.Range("G$r").Value = "System, backup@backdoor.com, system"
}

# Rows where "System, dnasr281@gmail.com" becomes "dnasr281@gmail.com, System"
$rowsDnasr = @(3, 6, 10, 11, 12, 13, 14, 15, 17, 18, 19, 20, 21, 22, 24, 30, 33, 37, 38, 39, 40, 41, 42, 44, 45, 46, 47, 48, 49, 51, 57, 60, 64, 65, 66, 67, 68, 69, 71, 72, 73, 74, 75, 76, 78, 86, 87, 88, 89, 93, 95, 96, 97, 99, 102, 112, 113, 114, 115, 119, 121, 122, 123, 125, 128, 138, 139, 140, 141, 145, 147, 148, 149, 151, 154)
foreach ($r in $rowsDnasr) {
    $ws.Range("G$r").Value = "dnasr281@gmail.com, System"
}

# Rows where "System, backup@backdoor.com" becomes "backup@backdoor.com, System"
$rowsBackup = @(4, 31, 58)
foreach ($r in $rowsBackup) {
    $ws.Range("G$r").Value = "backup@backdoor.com, System"
}

# Rows where "System, admin@admin.com" becomes "admin@admin.com, System"
$rowsAdmin = @(7, 34, 61)
foreach ($r in $rowsAdmin) {
    $ws.Range("G$r").Value = "admin@admin.com, System"
}
